$d = $word.ActiveDocument

# Update the date paragraph
$dateRange = $d.Paragraphs(1).Range
$dateRange.Text = "2024-02-08 Thursday"

# Update the table cells (5 columns x 20 rows), row-major order
$newValues = @(
    "99-68=",
    "21+63=",
    "26-15=",
    "0+38=",
    "61-56=",
    "52-38=",
    "78-26=",
    "43+38=",
    "45+7=",
    "15+13=",
    "3+23=",
    "97-85=",
    "84-14=",
    "24-6=",
    "51+8=",
    "34-21=",
    "52+42=",
    "42+13=",
    "97-15=",
    "92-59=",
    "98-91=",
    "69+26=",
    "83-33=",
    "61-25=",
    "60+19=",
    "98+0=",
    "51-9=",
    "51-49=",
    "69-48=",
    "87-16=",
    "37-36=",
    "82+6=",
    "57-26=",
    "92-6=",
    "78-54=",
    "55-17=",
    "97-30=",
    "59+12=",
    "52+31=",
    "15+40=",
    "48+50=",
    "94-79=",
    "50-12=",
    "59-14=",
    "36+51=",
    "15-5=",
    "18+78=",
    "1+51=",
    "82-19=",
    "72-12=",
    "81+15=",
    "18+8=",
    "67-37=",
    "56-8=",
    "95-21=",
    "27+19=",
    "97-18=",
    "97-90=",
    "83-17=",
    "19+27=",
    "37+59=",
    "31-19=",
    "58-7=",
    "24-21=",
    "47+23=",
    "41+27=",
    "60-34=",
    "0+63=",
    "90-76=",
    "21-15=",
    "47+5=",
    "68-29=",
    "56+38=",
    "64-46=",
    "50-32=",
    "30+62=",
    "70+8=",
    "75-36=",
    "70-11=",
    "42-6=",
    "4+32=",
    "4+14=",
    "9-5=",
    "69-45=",
    "98-59=",
    "86-73=",
    "82-28=",
    "47+47=",
    "85-11=",
    "97-64=",
    "54-1=",
    "40+30=",
    "59+4=",
    "40+25=",
    "45+29=",
    "62-47=",
    "33+63=",
    "30+5=",
    "23+9=",
    "88-73="
)

$t = $d.Tables(1)
$numCols = 5
$numRows = 20
$idx = 0
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $cellRange = $t.Cell($r, $c).Range
        $cellRange.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Done. Updated $idx cells."
